# Applies the workbook edit described by the diff:
#  1. Insert a new sheet "Đơn thu nợ" between "Đơn sale chính" and "Lương",
#     populated with a debt-collection report (headers + 1 data row + totals row).
#  2. In the "Lương" sheet, insert a new "Chiết khấu thu nợ tại <cơ sở>" row for each
#     location and add "Tổng lương tại <cơ sở>" rows before the grand total row,
#     shifting/renumbering the existing rows accordingly, and refresh totals.

$wb = $excel.ActiveWorkbook

# --- Step 1: insert "Đơn thu nợ" sheet right after "Đơn sale chính" ---
$wsSale = $wb.Worksheets.Item(1)
$wsDebt = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSale)
$wsDebt.Name = "Đơn thu nợ"

# Re-fetch the "Lương" sheet AFTER inserting the new sheet, since Worksheets.Item(n)
# is a positional reference and "Lương" has shifted from position 2 to position 3.
$wsLuong = $wb.Worksheets.Item(3)

# Header row (row 1, columns A:Y)
$wsDebt.Cells.Item(1, 1).Value = "Tiền tố"
$wsDebt.Cells.Item(1, 2).Value = "Mã đơn thu nợ"
$wsDebt.Cells.Item(1, 3).Value = "Lượng thu"
$wsDebt.Cells.Item(1, 4).Value = "Ngày thu"
$wsDebt.Cells.Item(1, 5).Value = "Cơ sở"
$wsDebt.Cells.Item(1, 6).Value = "Đơn nợ"
$wsDebt.Cells.Item(1, 7).Value = "Tên dịch vụ"
$wsDebt.Cells.Item(1, 8).Value = "Khách hàng"
$wsDebt.Cells.Item(1, 9).Value = "Nguồn khách"
$wsDebt.Cells.Item(1, 10).Value = "Sale chính"
$wsDebt.Cells.Item(1, 11).Value = "Đơn giá gốc"
$wsDebt.Cells.Item(1, 12).Value = "Sale phụ"
$wsDebt.Cells.Item(1, 13).Value = "Upsale"
$wsDebt.Cells.Item(1, 14).Value = "Đơn giá"
$wsDebt.Cells.Item(1, 15).Value = "Đã thanh toán"
$wsDebt.Cells.Item(1, 16).Value = "Bác sĩ 1"
$wsDebt.Cells.Item(1, 17).Value = "Bác sĩ 2"
$wsDebt.Cells.Item(1, 18).Value = "Tỉ lệ chiết khấu sale chính"
$wsDebt.Cells.Item(1, 19).Value = "Chiết khấu sale chính"
$wsDebt.Cells.Item(1, 20).Value = "Tỉ lệ chiết khấu sale phụ"
$wsDebt.Cells.Item(1, 21).Value = "Chiết khấu sale phụ"
$wsDebt.Cells.Item(1, 22).Value = "Tỉ lệ chiết khấu bác sĩ 1"
$wsDebt.Cells.Item(1, 23).Value = "Chiết khấu bác sĩ 1"
$wsDebt.Cells.Item(1, 24).Value = "Tỉ lệ chiết khấu bác sĩ 2"
$wsDebt.Cells.Item(1, 25).Value = "Chiết khấu bác sĩ 2"

# Keep the "Ngày thu" (Date collected) column as literal text, not an auto-parsed date
$wsDebt.Cells.Item(2, 4).NumberFormat = "@"

# Data row (row 2) - one collected-debt record
$wsDebt.Cells.Item(2, 1).Value = "TN"
$wsDebt.Cells.Item(2, 2).Value = 176
$wsDebt.Cells.Item(2, 3).Value = 500000
$wsDebt.Cells.Item(2, 4).Value = "08-04-2024"
$wsDebt.Cells.Item(2, 5).Value = "SÓC TRĂNG"
$wsDebt.Cells.Item(2, 6).Value = "HD-LUXURY-428"
$wsDebt.Cells.Item(2, 7).Value = "Phun môi"
$wsDebt.Cells.Item(2, 8).Value = "Lâm t.dân"
$wsDebt.Cells.Item(2, 9).Value = "Cá nhân"
$wsDebt.Cells.Item(2, 10).Value = "Cô Siêng giúp Việc"
$wsDebt.Cells.Item(2, 11).Value = 2500000
$wsDebt.Cells.Item(2, 12).Value = $null
$wsDebt.Cells.Item(2, 13).Value = $null
$wsDebt.Cells.Item(2, 14).Value = 2500000
$wsDebt.Cells.Item(2, 15).Value = 2500000
$wsDebt.Cells.Item(2, 16).Value = "CTV Ngoài"
$wsDebt.Cells.Item(2, 17).Value = $null
$wsDebt.Cells.Item(2, 18).Value = 0.1
$wsDebt.Cells.Item(2, 19).Value = 50000
$wsDebt.Cells.Item(2, 20).Value = 0
$wsDebt.Cells.Item(2, 21).Value = 0
$wsDebt.Cells.Item(2, 22).Value = 0
$wsDebt.Cells.Item(2, 23).Value = 0
$wsDebt.Cells.Item(2, 24).Value = 0
$wsDebt.Cells.Item(2, 25).Value = 0

# Totals row (row 3)
$wsDebt.Cells.Item(3, 1).Value = "Tổng"
$wsDebt.Cells.Item(3, 2).Value = 1
$wsDebt.Cells.Item(3, 3).Value = 500000
$wsDebt.Cells.Item(3, 4).Value = $null
$wsDebt.Cells.Item(3, 5).Value = $null
$wsDebt.Cells.Item(3, 6).Value = $null
$wsDebt.Cells.Item(3, 7).Value = $null
$wsDebt.Cells.Item(3, 8).Value = $null
$wsDebt.Cells.Item(3, 9).Value = $null
$wsDebt.Cells.Item(3, 10).Value = $null
$wsDebt.Cells.Item(3, 11).Value = 2500000
$wsDebt.Cells.Item(3, 12).Value = $null
$wsDebt.Cells.Item(3, 13).Value = 0
$wsDebt.Cells.Item(3, 14).Value = 2500000
$wsDebt.Cells.Item(3, 15).Value = 2500000
$wsDebt.Cells.Item(3, 16).Value = $null
$wsDebt.Cells.Item(3, 17).Value = $null
$wsDebt.Cells.Item(3, 18).Value = 0
$wsDebt.Cells.Item(3, 19).Value = 50000
$wsDebt.Cells.Item(3, 20).Value = 0
$wsDebt.Cells.Item(3, 21).Value = 0
$wsDebt.Cells.Item(3, 22).Value = 0
$wsDebt.Cells.Item(3, 23).Value = 0
$wsDebt.Cells.Item(3, 24).Value = 0
$wsDebt.Cells.Item(3, 25).Value = 0

# --- Step 2: update "Lương" sheet rows 11-37 (new "Chiết khấu thu nợ" + totals rows) ---
$wsLuong.Cells.Item(11, 1).Value = "Chiết khấu thu nợ tại CẦN THƠ"
$wsLuong.Cells.Item(11, 2).Value = 0
$wsLuong.Cells.Item(12, 1).Value = "Ứng lương tại CẦN THƠ"
$wsLuong.Cells.Item(12, 2).Value = 0
$wsLuong.Cells.Item(13, 1).Value = "Tổng công tại LONG XUYÊN"
$wsLuong.Cells.Item(13, 2).Value = 0
$wsLuong.Cells.Item(14, 1).Value = "Lương công tác tại LONG XUYÊN"
$wsLuong.Cells.Item(14, 2).Value = 0
$wsLuong.Cells.Item(15, 1).Value = "Lương cơ bản tại LONG XUYÊN"
$wsLuong.Cells.Item(15, 2).Value = $null
$wsLuong.Cells.Item(16, 1).Value = "Chiết khấu sale chính tại LONG XUYÊN"
$wsLuong.Cells.Item(16, 2).Value = 0
$wsLuong.Cells.Item(17, 1).Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$wsLuong.Cells.Item(17, 2).Value = 0
$wsLuong.Cells.Item(18, 1).Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$wsLuong.Cells.Item(18, 2).Value = 0
$wsLuong.Cells.Item(19, 1).Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$wsLuong.Cells.Item(19, 2).Value = 0
$wsLuong.Cells.Item(20, 1).Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$wsLuong.Cells.Item(20, 2).Value = 0
$wsLuong.Cells.Item(21, 1).Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$wsLuong.Cells.Item(21, 2).Value = 0
$wsLuong.Cells.Item(22, 1).Value = "Chiết khấu thu nợ tại LONG XUYÊN"
$wsLuong.Cells.Item(22, 2).Value = 0
$wsLuong.Cells.Item(23, 1).Value = "Ứng lương tại LONG XUYÊN"
$wsLuong.Cells.Item(23, 2).Value = 0
$wsLuong.Cells.Item(24, 1).Value = "Tổng công tại SÓC TRĂNG"
$wsLuong.Cells.Item(24, 2).Value = 7
$wsLuong.Cells.Item(25, 1).Value = "Lương cơ bản tại SÓC TRĂNG"
$wsLuong.Cells.Item(25, 2).Value = 1250000
$wsLuong.Cells.Item(26, 1).Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$wsLuong.Cells.Item(26, 2).Value = 0
$wsLuong.Cells.Item(27, 1).Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$wsLuong.Cells.Item(27, 2).Value = 0
$wsLuong.Cells.Item(28, 1).Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$wsLuong.Cells.Item(28, 2).Value = 0
$wsLuong.Cells.Item(29, 1).Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$wsLuong.Cells.Item(29, 2).Value = 0
$wsLuong.Cells.Item(30, 1).Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$wsLuong.Cells.Item(30, 2).Value = 0
$wsLuong.Cells.Item(31, 1).Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$wsLuong.Cells.Item(31, 2).Value = 0
$wsLuong.Cells.Item(32, 1).Value = "Chiết khấu thu nợ tại SÓC TRĂNG"
$wsLuong.Cells.Item(32, 2).Value = 50000
$wsLuong.Cells.Item(33, 1).Value = "Ứng lương tại SÓC TRĂNG"
$wsLuong.Cells.Item(33, 2).Value = 0
$wsLuong.Cells.Item(34, 1).Value = "Tổng lương tại CẦN THƠ"
$wsLuong.Cells.Item(34, 2).Value = 0
$wsLuong.Cells.Item(35, 1).Value = "Tổng lương tại LONG XUYÊN"
$wsLuong.Cells.Item(35, 2).Value = 0
$wsLuong.Cells.Item(36, 1).Value = "Tổng lương tại SÓC TRĂNG"
$wsLuong.Cells.Item(36, 2).Value = 1300000
$wsLuong.Cells.Item(37, 1).Value = "Tổng lương tại HỆ THỐNG"
$wsLuong.Cells.Item(37, 2).Value = 1300000
